$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 40 (pushing every
# subsequent record down by one row, up through the former row 87
# which becomes row 88).
$ws.Rows.Item(40).Insert()

$ws.Cells.Item(40, 1).Value = 4
$ws.Cells.Item(40, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(40, 3).Value = "Los Lagos"
$ws.Cells.Item(40, 4).Value = 45195
$ws.Cells.Item(40, 5).Value = 10
$ws.Cells.Item(40, 6).Value = 100112012
$ws.Cells.Item(40, 7).Value = "Espinaca"
$ws.Cells.Item(40, 8).Value = "Sin especificar"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 35
$ws.Cells.Item(40, 11).Value = 14000
$ws.Cells.Item(40, 12).Value = 14000
$ws.Cells.Item(40, 13).Value = 14000
$ws.Cells.Item(40, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(40, 15).Value = "Región Metropolitana"
$ws.Cells.Item(40, 16).Value = 1400
$ws.Cells.Item(40, 17).Value = 10
$ws.Cells.Item(40, 18).Value = "Hortaliza"
